$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "Map Construction" key/legend table that lives in columns N:P (rows
# 9-14) had a new "Map Flags" row inserted right above the existing "NPC
# Count" row. That pushes every row below it down by one, and the
# "Note: See example..." text that used to sit on row 14 now sits on row 15.
# Column A:D (the Tiles/NPCs legend) is untouched by this - only N:P move.
#
# Capture the current (pre-insert) N:P values for rows 9-14 first, then
# rewrite them one row lower, and finally stamp the new "Map Flags" row into
# row 9.
# ---------------------------------------------------------------------------

$srcN = @{}
$srcO = @{}
$srcP = @{}
$srcPFormula = @{}
$srcPIsFormula = @{}

for ($r = 9; $r -le 14; $r++) {
    $srcN[$r] = $ws.Range("N$r").Value2
    $srcO[$r] = $ws.Range("O$r").Value2
    if ($ws.Range("P$r").HasFormula) {
        $srcPIsFormula[$r] = $true
        $srcPFormula[$r] = $ws.Range("P$r").Formula
    } else {
        $srcPIsFormula[$r] = $false
        $srcP[$r] = $ws.Range("P$r").Value2
    }
}

# Shift rows 9-14 down to 10-15 (columns N:P only)
for ($r = 14; $r -ge 9; $r--) {
    $dest = $r + 1

    if ($null -eq $srcN[$r]) {
        $ws.Range("N$dest").ClearContents()
    } else {
        $ws.Range("N$dest").Value = $srcN[$r]
    }

    if ($null -eq $srcO[$r]) {
        $ws.Range("O$dest").ClearContents()
    } else {
        $ws.Range("O$dest").Value = $srcO[$r]
    }

    if ($srcPIsFormula[$r]) {
        $ws.Range("P$dest").Formula = $srcPFormula[$r]
    } elseif ($null -eq $srcP[$r]) {
        $ws.Range("P$dest").ClearContents()
    } else {
        $ws.Range("P$dest").Value = $srcP[$r]
    }
}

# New row 9: "Map Flags"
$ws.Range("N9").Value = "Map Flags"
$ws.Range("O9").Value = 1
$ws.Range("P9").Value = "-"

# Row 14's N:P got shifted into row 15 above, so row 14 itself is now blank.
$ws.Range("N14").ClearContents()
$ws.Range("O14").ClearContents()
$ws.Range("P14").ClearContents()

# ---------------------------------------------------------------------------
# Selection moved to P9 (this is where the row was inserted), with the view
# scrolled so row 3 is at the top.
# ---------------------------------------------------------------------------
$ws.Range("P9").Select()
$excel.ActiveWindow.ScrollRow = 3

# ---------------------------------------------------------------------------
# The key-legend picture anchored below the table shifts down by one row
# along with the table (its anchor offsets change even though its size
# doesn't).
# ---------------------------------------------------------------------------
$pic = $ws.Shapes.Item(1)
$pic.Left = 741
$pic.Top = 239.25
$pic.Width = 216.75
$pic.Height = 201
